$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) Insert a new (underlined) heading paragraph before the existing first
#    paragraph ("ALTER TABLE usuario ...").
# ---------------------------------------------------------------------------
$firstPara = $d.Paragraphs(1)
$firstPara.Range.InsertParagraphBefore()

$headingXml = "<w:p $wns>" +
    "<w:pPr><w:rPr><w:u w:val=`"single`"/></w:rPr></w:pPr>" +
    "<w:r><w:t xml:space=`"preserve`">Alter para criar coluna </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>email</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`"> na tabela </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>Usuario</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "</w:p>"
$d.Paragraphs(1).Range.InsertXML($headingXml)

# ---------------------------------------------------------------------------
# 2) After the original two paragraphs (now paragraphs 2 and 3), append 8
#    more paragraphs: an empty one, the "Delete..." explanation, the
#    TRUNCATE statement, an empty one, the "Alter..." explanation, the new
#    ALTER TABLE statement, its ADD COLUMN statement, and a trailing empty
#    paragraph.
# ---------------------------------------------------------------------------
$anchor = $d.Paragraphs(3)
for ($i = 0; $i -lt 8; $i++) {
    $anchor.Range.InsertParagraphAfter()
}

# Paragraph 4: blank separator
$d.Paragraphs(4).Range.InsertXML("<w:p $wns/>")

# Paragraph 5: "Delete para adicionar os dados de tipoUsuario ..."
$deleteXml = "<w:p $wns>" +
    "<w:r><w:t xml:space=`"preserve`">Delete para adicionar os dados de </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>tipoUsuario</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`"> " + [char]0x2013 + " enviar esse comando ap" + [char]0x00F3 + "s enviar o </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>alter</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>table</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`"> acima</w:t></w:r>" +
    "</w:p>"
$d.Paragraphs(5).Range.InsertXML($deleteXml)

# Paragraph 6: "TRUNCATE TABLE usuario CASCADE;"
$truncateXml = "<w:p $wns>" +
    "<w:r><w:t xml:space=`"preserve`">TRUNCATE TABLE </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>usuario</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`"> CASCADE;</w:t></w:r>" +
    "</w:p>"
$d.Paragraphs(6).Range.InsertXML($truncateXml)

# Paragraph 7: blank separator
$d.Paragraphs(7).Range.InsertXML("<w:p $wns/>")

# Paragraph 8: "Alter para criar coluna tipoUsuario na tabela Usuario"
$alterHeadingXml = "<w:p $wns>" +
    "<w:r><w:t xml:space=`"preserve`">Alter para criar coluna </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>tipoUsuario</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`"> na tabela </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>Usuario</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "</w:p>"
$d.Paragraphs(8).Range.InsertXML($alterHeadingXml)

# Paragraph 9: "ALTER TABLE usuario"
$alterTableXml = "<w:p $wns>" +
    "<w:r><w:t xml:space=`"preserve`">ALTER TABLE </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>usuario</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "</w:p>"
$d.Paragraphs(9).Range.InsertXML($alterTableXml)

# Paragraph 10: "ADD COLUMN tipoUsuario VARCHAR(15);"
$addColumnXml = "<w:p $wns>" +
    "<w:r><w:t xml:space=`"preserve`">ADD COLUMN </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>tipoUsuario</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>" +
    "<w:proofErr w:type=`"gramStart`"/>" +
    "<w:r><w:t>VARCHAR(</w:t></w:r>" +
    "<w:proofErr w:type=`"gramEnd`"/>" +
    "<w:r><w:t>15);</w:t></w:r>" +
    "</w:p>"
$d.Paragraphs(10).Range.InsertXML($addColumnXml)

# Paragraph 11: trailing blank paragraph
$d.Paragraphs(11).Range.InsertXML("<w:p $wns/>")

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
